$d = $word.ActiveDocument

# Locate the paragraph that ends the "Solution b)" answer for Problem #2;
# the new Problem #3 block is inserted right after the blank paragraph
# that follows it.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*pair of white socks.*") {
        $anchorIndex = $i
        break
    }
}
$anchor = $d.Paragraphs.Item($anchorIndex + 1)

# Create one fresh blank paragraph right after the anchor, then replace its
# (placeholder) content with the real Problem #3 block via OOXML so the
# paragraph/run formatting matches exactly.
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$placeholder = $d.Paragraphs.Item($anchorIndex + 2)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Problem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> #3: &#8220;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>Predicting Fingers</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>&#8221;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:i/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:i/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>A little girl counts using the fingers of her left hand as follows: She starts by calling her thumb 1, the first finger 2, middle finder 3, ring finger 4, and little finger 5. Then she reverses direction, calling the ring finger 6, middle finger 7, first finger 8 and thumb 9, after which she calls her first finger 10 and so on. If she continues to count in this manner, on which finger will she stop? a) What if the girl counts from 1 to 10 b) What if the girl counts from 1 to 100 c) What if the girl counts from 1 to 1000</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:i/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:i/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
'@

$placeholder.Range.InsertXML($xml)
